# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# timestamps on the zh-cn and de-de report sheets (regenerated handback report).

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-22 06:56:49"
$wsZhCn.Range("H2").Value = "2016-03-22 06:57:34"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-22 06:56:57"
$wsDeDe.Range("H2").Value = "2016-03-22 06:57:48"
